$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.083.15"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.349.55"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'545.19"
$ws.Range("E5").Value = "  +5.80%  "
$ws.Range("D6").Value = "'134.94"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "2.348.44"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("E13").Value = "  +6.34%  "
$ws.Range("D14").Value = "2.767.29"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "'23.59"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "58.058.80"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "2.352.08"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'10.65"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "'334.90"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D24").Value = "'61.82"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'8.45"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "'1.41"
$ws.Range("E28").Value = "  +7.79%  "
$ws.Range("E29").Value = "  +5.38%  "
$ws.Range("D30").Value = "'170.14"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").Value = "'6.14"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +18.19%  "
$ws.Range("D34").Value = "'18.44"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'4.20"
$ws.Range("E37").Value = "  +6.60%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("D40").Value = "'39.23"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("D41").Value = "'147.50"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").Value = "'0.379"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'287.40"
$ws.Range("E43").Value = "  +4.34%  "
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").Value = "'19.25"
$ws.Range("E45").Value = "  +5.71%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "'0.562"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("E51").Value = "  +0.66%  "

# Remove the quote-prefix formatting artifact introduced by forcing
# numeric-looking price strings to stay as text, so cell styles are
# left exactly as they were before the edits.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D48").ClearFormats()
